# Update "want to go" counts (column F) and a couple of cover image links
# (column I) across the "展览" (sheet 1), "演出" (sheet 2) and "全部类型"
# (sheet 4) worksheets, to match freshly re-generated data.

$wb = $excel.ActiveWorkbook

# ----- Sheet "展览" (Exhibitions) -----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 14
$ws.Range("F3").Value = 968
$ws.Range("F6").Value = 1087
$ws.Range("F7").Value = 869
$ws.Range("F8").Value = 262
$ws.Range("F11").Value = 858
$ws.Range("F13").Value = 587
$ws.Range("F15").Value = 1351
$ws.Range("F18").Value = 1221
$ws.Range("F20").Value = 1481
$ws.Range("F21").Value = 725
$ws.Range("F22").Value = 205
$ws.Range("F23").Value = 1291
$ws.Range("I24").Value = "//i2.hdslb.com/bfs/openplatform/202404/KVZUuqc51713150203436.jpeg"
$ws.Range("F25").Value = 1040
$ws.Range("F26").Value = 364
$ws.Range("F27").Value = 3204
$ws.Range("F30").Value = 1429

# ----- Sheet "演出" (Performances) -----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 160
$ws.Range("F9").Value = 15

# ----- Sheet "全部类型" (All types) -----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 14
$ws.Range("F6").Value = 968
$ws.Range("F10").Value = 1087
$ws.Range("F11").Value = 869
$ws.Range("F12").Value = 262
$ws.Range("F17").Value = 160
$ws.Range("F20").Value = 15
$ws.Range("F21").Value = 858
$ws.Range("F23").Value = 587
$ws.Range("F25").Value = 1351
$ws.Range("F28").Value = 1221
$ws.Range("F30").Value = 1481
$ws.Range("F31").Value = 725
$ws.Range("F32").Value = 205
$ws.Range("F33").Value = 1291
$ws.Range("I34").Value = "//i2.hdslb.com/bfs/openplatform/202404/KVZUuqc51713150203436.jpeg"
$ws.Range("F37").Value = 1040
$ws.Range("F38").Value = 364
$ws.Range("F39").Value = 3204
$ws.Range("F42").Value = 1429
